# Update gh-pages "想去人数" (want-to-go count) figures for the 苏州-漫展信息
# workbook. The same set of events is listed on the "展览" sheet and again on
# the "全部类型" sheet, so every changed event needs its F-column value
# updated in both places.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row number on "展览", row number on "全部类型", new value for column F.
$updates = @(
    @(3,  3,  506),
    @(4,  4,  1490),
    @(5,  5,  146),
    @(6,  6,  138),
    @(9,  10, 132),
    @(10, 11, 723),
    @(11, 13, 1038),
    @(12, 14, 59),
    @(13, 15, 311),
    @(14, 16, 44),
    @(15, 18, 6315),
    @(16, 19, 82),
    @(20, 23, 15141),
    @(21, 24, 1501),
    @(22, 25, 268),
    @(24, 27, 97),
    @(25, 28, 10976),
    @(26, 29, 727),
    @(27, 30, 4284),
    @(28, 31, 220),
    @(30, 33, 8)
)

foreach ($u in $updates) {
    $row1 = $u[0]
    $row4 = $u[1]
    $newVal = $u[2]

    $ws1.Cells.Item($row1, 6).Value = $newVal
    $ws4.Cells.Item($row4, 6).Value = $newVal
}
